$wb = $excel.ActiveWorkbook

$wsReview = $wb.Worksheets.Item("LH_WF_PUBLISHAUDIO_REVIEW")
$wsHistory = $wb.Worksheets.Item("VERSION-HISTORY")

# --- Update the "Owner Status" column (H) on the review sheet ---
# Row 2 and Row 4 and Row 5 move from "open" to "not applicable"
# Row 3 moves from "open" to "closed"
$wsReview.Range("H2").Value = "not applicable"
$wsReview.Range("H3").Value = "closed"
$wsReview.Range("H4").Value = "not applicable"
$wsReview.Range("H5").Value = "not applicable"

# --- Add a new entry to the VERSION-HISTORY sheet describing this update ---
$wsHistory.Range("A3").Value = "v1.1"
$wsHistory.Range("B3").Value = "eman"
$wsHistory.Range("C3").Value = "edit the owner status of the sheet"
$wsHistory.Range("D3").Value = (Get-Date -Year 2025 -Month 4 -Day 29 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

# --- Update view state: selection moves, and VERSION-HISTORY becomes the active tab ---
$wsReview.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$wsReview.Range("H5").Select()

$wsHistory.Activate()
$wsHistory.Range("D7").Select()
